$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Full1")

# --- Header for new column F ---
$ws.Range("F1").Value = "significance"

# --- Apply the existing plain-font style (already used by column D/E from row 50 on)
# to the new/updated F column cells, by copying format from a cell that already uses it ---
$ws.Range("D50").Copy()
$ws.Range("F2:F75").PasteSpecial(-4122)
$ws.Range("D50").Copy()
$ws.Range("F76:F77").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Fill in significance yes/no values for rows 2-75 ---
$ws.Cells.Item(2, 6).Value = "yes"
$ws.Cells.Item(3, 6).Value = "no"
$ws.Cells.Item(4, 6).Value = "yes"
$ws.Cells.Item(5, 6).Value = "no"
$ws.Cells.Item(6, 6).Value = "yes"
$ws.Cells.Item(7, 6).Value = "yes"
$ws.Cells.Item(8, 6).Value = "yes"
$ws.Cells.Item(9, 6).Value = "yes"
$ws.Cells.Item(10, 6).Value = "no"
$ws.Cells.Item(11, 6).Value = "no"
$ws.Cells.Item(12, 6).Value = "no"
$ws.Cells.Item(13, 6).Value = "no"
$ws.Cells.Item(14, 6).Value = "no"
$ws.Cells.Item(15, 6).Value = "no"
$ws.Cells.Item(16, 6).Value = "yes"
$ws.Cells.Item(17, 6).Value = "yes"
$ws.Cells.Item(18, 6).Value = "yes"
$ws.Cells.Item(19, 6).Value = "yes"
$ws.Cells.Item(20, 6).Value = "no"
$ws.Cells.Item(21, 6).Value = "no"
$ws.Cells.Item(22, 6).Value = "yes"
$ws.Cells.Item(23, 6).Value = "yes"
$ws.Cells.Item(24, 6).Value = "yes"
$ws.Cells.Item(25, 6).Value = "yes"
$ws.Cells.Item(26, 6).Value = "yes"
$ws.Cells.Item(27, 6).Value = "yes"
$ws.Cells.Item(28, 6).Value = "yes"
$ws.Cells.Item(29, 6).Value = "yes"
$ws.Cells.Item(30, 6).Value = "yes"
$ws.Cells.Item(31, 6).Value = "yes"
$ws.Cells.Item(32, 6).Value = "yes"
$ws.Cells.Item(33, 6).Value = "yes"
$ws.Cells.Item(34, 6).Value = "yes"
$ws.Cells.Item(35, 6).Value = "yes"
$ws.Cells.Item(36, 6).Value = "yes"
$ws.Cells.Item(37, 6).Value = "yes"
$ws.Cells.Item(38, 6).Value = "no"
$ws.Cells.Item(39, 6).Value = "no"
$ws.Cells.Item(40, 6).Value = "no"
$ws.Cells.Item(41, 6).Value = "yes"
$ws.Cells.Item(42, 6).Value = "no"
$ws.Cells.Item(43, 6).Value = "yes"
$ws.Cells.Item(44, 6).Value = "yes"
$ws.Cells.Item(45, 6).Value = "yes"
$ws.Cells.Item(46, 6).Value = "yes"
$ws.Cells.Item(47, 6).Value = "yes"
$ws.Cells.Item(48, 6).Value = "yes"
$ws.Cells.Item(49, 6).Value = "yes"
$ws.Cells.Item(50, 6).Value = "yes"
$ws.Cells.Item(51, 6).Value = "yes"
$ws.Cells.Item(52, 6).Value = "no"
$ws.Cells.Item(53, 6).Value = "yes"
$ws.Cells.Item(54, 6).Value = "no"
$ws.Cells.Item(55, 6).Value = "yes"
$ws.Cells.Item(56, 6).Value = "yes"
$ws.Cells.Item(57, 6).Value = "yes"
$ws.Cells.Item(58, 6).Value = "yes"
$ws.Cells.Item(59, 6).Value = "yes"
$ws.Cells.Item(60, 6).Value = "no"
$ws.Cells.Item(61, 6).Value = "yes"
$ws.Cells.Item(62, 6).Value = "yes"
$ws.Cells.Item(63, 6).Value = "no"
$ws.Cells.Item(64, 6).Value = "no"
$ws.Cells.Item(65, 6).Value = "yes"
$ws.Cells.Item(66, 6).Value = "no"
$ws.Cells.Item(67, 6).Value = "no"
$ws.Cells.Item(68, 6).Value = "no"
$ws.Cells.Item(69, 6).Value = "no"
$ws.Cells.Item(70, 6).Value = "no"
$ws.Cells.Item(71, 6).Value = "no"
$ws.Cells.Item(72, 6).Value = "yes"
$ws.Cells.Item(73, 6).Value = "no"
$ws.Cells.Item(74, 6).Value = "no"
$ws.Cells.Item(75, 6).Value = "yes"

# --- Restore view/selection state ---
$ws.Range("F56").Select()
$excel.ActiveWindow.ScrollRow = 51

